# Applies the "Add files via upload" edit to Sentiment_Analysis.pptx / slide 1:
#  1) Name text box ("Tejeshwar M" -> "santhosh s"), run language en-IN -> en-US
#  2) Roll number text box ("813821104109" -> "813821104087")
#  3) Department text box ("CSE B" -> "COMPUTER SCIENCE AND ENGINEERING"),
#     run language en-IN -> en-US, and a resize/reposition of the box.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) "Tejeshwar M" -> "santhosh s" -------------------------------------
$nameShape = $s.Shapes.Item(4)
$nameRange = $nameShape.TextFrame.TextRange
$nameRange.Text = "santhosh s"
$nameRange.LanguageID = "en-US"

# --- 2) Roll number -------------------------------------------------------
$rollShape = $s.Shapes.Item(7)
$rollShape.TextFrame.TextRange.Text = "813821104087"

# --- 3) Department text box: text, language, position & size -------------
$deptShape = $s.Shapes.Item(8)

# Floats below are chosen (in points) so that this host's pt -> EMU
# conversion reproduces the exact target EMU offsets/extents:
#   off  x="6465958" y="3503014"
#   ext  cx="3211441" cy="646331"
$deptShape.Left = 509.13055118110236
$deptShape.Top = 275.827874015748
$deptShape.Width = 252.86938007874016
$deptShape.Height = 50.89221472440945

$deptRange = $deptShape.TextFrame.TextRange
$deptRange.Text = "COMPUTER SCIENCE AND ENGINEERING"
$deptRange.LanguageID = "en-US"
